# Updates the cryptocurrency price/volume snapshot in columns D (Price) and
# E (Volume(1h)) for rows 2-51, per the Thu Jun 8 17:27:46 UTC 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value. Values that look like plain numbers (e.g. "263.33")
# get a leading apostrophe so Excel stores them as text, matching the
# original inline-string cells (prices such as "26.523.13" already contain
# two dots and are never number-like, so they need no marker).
$updates = [ordered]@{
    "D2" = "26.523.13"
    "E2" = "  +0.42%  "
    "D3" = "1.847.23"
    "E3" = "  +0.31%  "
    "E4" = "  +0.02%  "
    "D5" = "`'263.33"
    "E5" = "  +1.13%  "
    "E6" = "  +0.03%  "
    "D7" = "`'0.5224"
    "E7" = "  +1.83%  "
    "D8" = "`'0.3234"
    "E8" = "  +0.86%  "
    "E9" = "  +0.23%  "
    "D10" = "`'18.67"
    "E10" = "  -0.17%  "
    "D11" = "`'0.7711"
    "E11" = "  +0.35%  "
    "D12" = "`'0.07775"
    "E12" = "  +1.24%  "
    "D13" = "1.865.20"
    "E13" = "  -0.42%  "
    "D14" = "`'88.33"
    "E14" = "  -0.32%  "
    "D15" = "`'5.006"
    "E15" = "  -0.22%  "
    "E16" = "  +0.04%  "
    "D17" = "`'13.92"
    "E17" = "  -0.93%  "
    "E18" = "  -0.01%  "
    "D19" = "`'0.000007923"
    "E19" = "  +0.46%  "
    "D20" = "26.554.84"
    "E20" = "  +0.32%  "
    "D21" = "2.090.36"
    "E21" = "  -0.51%  "
    "D22" = "`'4.613"
    "E22" = "  +1.00%  "
    "D23" = "`'9.424"
    "E23" = "  -1.17%  "
    "D24" = "`'5.957"
    "E24" = "  +0.14%  "
    "D25" = "`'142.99"
    "E25" = "  -1.23%  "
    "D26" = "`'2.180"
    "E26" = "  -6.14%  "
    "E27" = "  +0.52%  "
    "E28" = "  +0.14%  "
    "D29" = "`'111.75"
    "E29" = "  +0.70%  "
    "D30" = "`'4.164"
    "E30" = "  -0.11%  "
    "D31" = "`'0.08739"
    "E31" = "  +0.35%  "
    "D32" = "`'4.109"
    "E32" = "  -1.05%  "
    "D33" = "`'0.04823"
    "E33" = "  +0.18%  "
    "E34" = "  -0.29%  "
    "D35" = "`'2.871"
    "E35" = "  +1.09%  "
    "D36" = "`'0.7147"
    "E36" = "  +4.37%  "
    "D38" = "`'0.01784"
    "E38" = "  -1.18%  "
    "D39" = "`'2.182"
    "E39" = "  -0.94%  "
    "D40" = "`'0.4832"
    "E40" = "  -1.67%  "
    "D41" = "`'112.26"
    "E41" = "  -0.89%  "
    "D42" = "`'0.8955"
    "E42" = "  -0.82%  "
    "D43" = "`'6.043"
    "E43" = "  -1.38%  "
    "E44" = "  +0.04%  "
    "D45" = "`'7.620"
    "E45" = "  -1.65%  "
    "D46" = "`'0.4161"
    "E46" = "  -1.93%  "
    "D47" = "`'0.05898"
    "E47" = "  +0.14%  "
    "D48" = "`'9.024"
    "E48" = "  -1.28%  "
    "D49" = "`'34.89"
    "E49" = "  -0.12%  "
    "D50" = "`'0.1225"
    "D51" = "`'0.8842"
    "E51" = "  +4.05%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
